# feat: add 2022-Q3 data
#
# - insert a new "2022-Q3" sheet right after "总计" (pushing the existing
#   "2022-Q2" / "2022-Q1" sheets one slot later)
# - populate it with the Q3 fund-holding table
# - prepend a 2022-Q3 summary row to the "总计" sheet, shifting the
#   existing 2022-Q2 / 2022-Q1 rows down by one

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)          # "总计"

# NOTE: worksheet handles returned by Worksheets.Item(...) track a *slot*,
# not a stable identity in this host — inserting a sheet before one shifts
# what that handle resolves to. So we deliberately re-fetch "2022-Q1" by
# name at the very end, after all sheet insertions are already done.

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet right after "总计"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $total)
$q3.Name = "2022-Q3"

# Header row (B1:H1)
$headers = New-Object 'object[,]' 1,7
$headers[0,0] = "基金代码"
$headers[0,1] = "基金名称"
$headers[0,2] = "基金规模"
$headers[0,3] = "股票总仓位"
$headers[0,4] = "仓位占比"
$headers[0,5] = "持有市值(亿元)"
$headers[0,6] = "仓位排名"
$q3.Range("B1:H1").Value = $headers

# Row index column (A2:A12) — 0-based row counter
$idx = New-Object 'object[,]' 11,1
for ($i = 0; $i -lt 11; $i++) { $idx[$i,0] = $i }
$q3.Range("A2:A12").Value = $idx

# Fund code / name / size / position / ratio / market-value / rank
$codes = New-Object 'object[,]' 11,1
$codes[0,0]  = "009714"
$codes[1,0]  = "008371"
$codes[2,0]  = "040016"
$codes[3,0]  = "040020"
$codes[4,0]  = "001825"
$codes[5,0]  = "014271"
$codes[6,0]  = "015071"
$codes[7,0]  = "014380"
$codes[8,0]  = "014272"
$codes[9,0]  = "015072"
$codes[10,0] = "014976"

$names = New-Object 'object[,]' 11,1
$names[0,0]  = "华安聚优精选混合"
$names[1,0]  = "华安汇智精选混合"
$names[2,0]  = "华安行业轮动混合"
$names[3,0]  = "华安升级主题混合A"
$names[4,0]  = "建信中国制造2025股票A"
$names[5,0]  = "大成北交所两年定开混合A"
$names[6,0]  = "鑫元专精特新混合A"
$names[7,0]  = "建信中国制造2025股票C"
$names[8,0]  = "大成北交所两年定开混合C"
$names[9,0]  = "鑫元专精特新混合C"
$names[10,0] = "华安升级主题混合C"
$q3.Range("C2:C12").Value = $names

$sizes = New-Object 'object[,]' 11,1
$sizes[0,0]  = "66.05"
$sizes[1,0]  = "7.01"
$sizes[2,0]  = "6.75"
$sizes[3,0]  = "4.40"
$sizes[4,0]  = "3.96"
$sizes[5,0]  = "3.45"
$sizes[6,0]  = "2.65"
$sizes[7,0]  = "1.69"
$sizes[8,0]  = "0.82"
$sizes[9,0]  = "0.25"
$sizes[10,0] = "0.00"

$positions = New-Object 'object[,]' 11,1
$positions[0,0]  = "90.02"
$positions[1,0]  = "91.18"
$positions[2,0]  = "79.22"
$positions[3,0]  = "85.97"
$positions[4,0]  = "86.05"
$positions[5,0]  = "65.31"
$positions[6,0]  = "74.01"
$positions[7,0]  = "86.05"
$positions[8,0]  = "65.31"
$positions[9,0]  = "74.01"
$positions[10,0] = "85.97"

$ratios = New-Object 'object[,]' 11,1
$ratios[0,0]  = "2.58"
$ratios[1,0]  = "2.60"
$ratios[2,0]  = "2.46"
$ratios[3,0]  = "2.69"
$ratios[4,0]  = "2.94"
$ratios[5,0]  = "2.75"
$ratios[6,0]  = "2.64"
$ratios[7,0]  = "2.94"
$ratios[8,0]  = "2.75"
$ratios[9,0]  = "2.64"
$ratios[10,0] = "2.69"

$mktval = New-Object 'object[,]' 10,1
$mktval[0,0] = "1.7041"
$mktval[1,0] = "0.1823"
$mktval[2,0] = "0.1660"
$mktval[3,0] = "0.1184"
$mktval[4,0] = "0.1164"
$mktval[5,0] = "0.0949"
$mktval[6,0] = "0.0700"
$mktval[7,0] = "0.0497"
$mktval[8,0] = "0.0226"
$mktval[9,0] = "0.0066"
# row 12's market value is the literal number 0 (not text) — set separately

$ranks = New-Object 'object[,]' 11,1
$ranks[0,0]  = 10
$ranks[1,0]  = 9
$ranks[2,0]  = 9
$ranks[3,0]  = 10
$ranks[4,0]  = 10
$ranks[5,0]  = 8
$ranks[6,0]  = 4
$ranks[7,0]  = 10
$ranks[8,0]  = 8
$ranks[9,0]  = 4
$ranks[10,0] = 10

# Write the text-valued columns via a NumberFormat="@" + PasteSpecial(values-only)
# round-trip so they land as genuine text ("66.05") instead of being
# auto-coerced to numbers by the COM value setter, while keeping the
# default (unstyled) cell format.
$blank = $q3.Range("Z1")
foreach ($colSpec in @(
        @{ Range = "B2:B12"; Data = $codes },
        @{ Range = "D2:D12"; Data = $sizes },
        @{ Range = "E2:E12"; Data = $positions },
        @{ Range = "F2:F12"; Data = $ratios },
        @{ Range = "G2:G11"; Data = $mktval }
    )) {
    $target = $q3.Range($colSpec.Range)
    $target.NumberFormat = "@"
    $target.Value = $colSpec.Data
    $blank.Copy()
    $target.PasteSpecial(-4122)
}
$blank.Clear()

# G12 is a genuine number 0 (not text)
$q3.Range("G12").Value = 0

$q3.Range("H2:H12").Value = $ranks

# ---------------------------------------------------------------------
# Formatting: replicate the "s=2" bold/border style used by the header
# row and the A-column index cells elsewhere in the workbook, by
# copy/pasting formats from the already-styled "总计" sheet.
# ---------------------------------------------------------------------
$total.Range("B1:D1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$total.Range("A2").Copy()
$q3.Range("A2:A12").PasteSpecial(-4122)

$q3.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. Prepend a 2022-Q3 row to the "总计" summary sheet, shifting the
#    existing 2022-Q2 / 2022-Q1 rows down by one.
#
#    We avoid Rows.Insert() here (it interpolates a new, unwanted cell
#    style for the shifted-into blank row) and instead shift the data
#    down manually, bottom row first, then overwrite row 2 with the new
#    2022-Q3 figures.
# ---------------------------------------------------------------------
$oldA2 = $total.Cells.Item(2, 1).Value()
$oldB2 = $total.Cells.Item(2, 2).Value()
$oldC2 = $total.Cells.Item(2, 3).Value()
$oldD2 = $total.Cells.Item(2, 4).Value()
$oldA3 = $total.Cells.Item(3, 1).Value()
$oldB3 = $total.Cells.Item(3, 2).Value()
$oldC3 = $total.Cells.Item(3, 3).Value()
$oldD3 = $total.Cells.Item(3, 4).Value()

# extend the A-column "s=2" style down onto the new row 4
$total.Range("A2").Copy()
$total.Range("A2:A4").PasteSpecial(-4122)

# row 4 <- old row 3 ("2022-Q1"), index bumped 1 -> 2
$total.Cells.Item(4, 1).Value = $oldA3 + 1
$total.Cells.Item(4, 2).Value = $oldB3
$total.Cells.Item(4, 3).Value = $oldC3
$total.Cells.Item(4, 4).Value = $oldD3

# row 3 <- old row 2 ("2022-Q2"), index bumped 0 -> 1
$total.Cells.Item(3, 1).Value = $oldA2 + 1
$total.Cells.Item(3, 2).Value = $oldB2
$total.Cells.Item(3, 3).Value = $oldC2
$total.Cells.Item(3, 4).Value = $oldD2

# row 2 <- new "2022-Q3" figures
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 11
$total.Cells.Item(2, 4).Value = 2.53

$total.Range("A1").Select()

# ---------------------------------------------------------------------
# Restore the originally-active sheet ("2022-Q1") as the selected tab,
# since adding the new sheet above moved the COM selection. Re-fetched
# by name now that the sheet collection is in its final shape.
# ---------------------------------------------------------------------
$oldQ1 = $wb.Worksheets.Item("2022-Q1")
$oldQ1.Activate()
$oldQ1.Range("A1").Select()
